# section_list_example.xlsx - header row rework
#
# The sheet's header row (row 1) is replaced: several columns are dropped
# (title, credits, dept_name, time, day, instructor_id), section_id moves
# next to course_id, and two new columns (start, end) plus a renamed
# classroom/lesson/limit block and a final "dat" column take their place.
# The end result only spans A1:H1 (down from A1:K1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K (11) carried a leftover custom width (15.6640625) that isn't
# tied to any header text. Inserting a blank column at K pushes that
# width formatting out to column L (12), matching how the workbook ended
# up with <col min="12" .../> instead of <col min="11" .../>.
$ws.Columns.Item(11).Insert()

# Clear the whole old header row (now A1:L1 after the insert above) and
# write the new, shorter header row in its place.
$ws.Range("A1:L1").ClearContents()

$ws.Range("A1").Value = "course_id"
$ws.Range("B1").Value = "section_id"
$ws.Range("C1").Value = "start"
$ws.Range("D1").Value = "end"
$ws.Range("E1").Value = "classroom_no"
$ws.Range("F1").Value = "lesson"
$ws.Range("G1").Value = "limit"
$ws.Range("H1").Value = "dat"

# Match the saved selection/active cell from the edit.
$ws.Range("F5").Select()
